# Apply the commit's changes:
#  1. Rename the worksheet from "AlphaFiberF-HW10.xpc" to "AlphaFiberF"
#  2. Tiny floating-point refinements on row 13 (C13, F13, L13, M13) - result
#     of re-exporting with the new Gaussian-Quadrature routine
#  3. Append a new data row (row 16) for the "HexGrid-60degTilt5degRes" scheme

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet/tab
$ws.Name = "AlphaFiberF"

# 2. Refresh the slightly-changed values in row 13
$ws.Range("C13").Value = 0.9818658657972137
$ws.Range("F13").Value = 0.9818658657972137
$ws.Range("L13").Value = 0.9894838732174307
$ws.Range("M13").Value = 0.9926044986113153

# 3. Append new row 16 with the HexGrid-60degTilt5degRes data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.732273679167148
$ws.Range("D16").Value = 1.686442167232839
$ws.Range("E16").Value = 0.7414874669677609
$ws.Range("F16").Value = 1.732273679167148
$ws.Range("G16").Value = 0.9880047329540944
$ws.Range("H16").Value = 1.071215803390088
$ws.Range("I16").Value = 0.8604059345742562
$ws.Range("J16").Value = 1.686442167232839
$ws.Range("K16").Value = 1.2139648171003
$ws.Range("L16").Value = 1.473119248133724
$ws.Range("M16").Value = 1.179971630714365

# Match the styling used by the rest of column A (bold/border/center "header" style)
# by copying the format from the cell directly above (A15) rather than
# rebuilding it property-by-property.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
